$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "26.338.62"
Set-TextCell "E2" "  +0.16%  "

# Row 3
Set-TextCell "D3" "1.685.71"
Set-TextCell "E3" "  +0.90%  "

# Row 4
Set-TextCell "D4" "1.009"
Set-TextCell "E4" "  +0.42%  "

# Row 5
Set-TextCell "D5" "218.31"
Set-TextCell "E5" "  +0.24%  "

# Row 6
Set-TextCell "D6" "0.5235"
Set-TextCell "E6" "  +3.03%  "

# Row 7
Set-TextCell "E7" "  +0.42%  "

# Row 8
Set-TextCell "D8" "0.2700"
Set-TextCell "E8" "  +1.61%  "

# Row 9
Set-TextCell "B9" "Dogecoin"
Set-TextCell "C9" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell "D9" "0.06406"
Set-TextCell "E9" "  +0.86%  "

# Row 10
Set-TextCell "B10" "Solana"
Set-TextCell "C10" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell "D10" "22.00"
Set-TextCell "E10" "  +2.12%  "

# Row 11
Set-TextCell "E11" "  +1.79%  "

# Row 12
Set-TextCell "D12" "1.721.64"
Set-TextCell "E12" "  +3.20%  "

# Row 13
Set-TextCell "D13" "4.560"
Set-TextCell "E13" "  +0.19%  "

# Row 14
Set-TextCell "D14" "0.5788"
Set-TextCell "E14" "  -0.49%  "

# Row 15
Set-TextCell "D15" "0.000008475"
Set-TextCell "E15" "  -0.56%  "

# Row 16
Set-TextCell "D16" "64.26"
Set-TextCell "E16" "  -0.79%  "

# Row 17
Set-TextCell "D17" "26.388.12"
Set-TextCell "E17" "  +0.12%  "

# Row 18
Set-TextCell "D18" "4.919"
Set-TextCell "E18" "  -0.41%  "

# Row 19
Set-TextCell "E19" "  +0.32%  "

# Row 20
Set-TextCell "D20" "10.87"
Set-TextCell "E20" "  +0.22%  "

# Row 21
Set-TextCell "D21" "188.56"
Set-TextCell "E21" "  -0.14%  "

# Row 22
Set-TextCell "D22" "6.189"
Set-TextCell "E22" "  -0.22%  "

# Row 23
Set-TextCell "D23" "1.010"
Set-TextCell "E23" "  +0.41%  "

# Row 24
Set-TextCell "D24" "144.44"
Set-TextCell "E24" "  +0.77%  "

# Row 25
Set-TextCell "D25" "7.708"
Set-TextCell "E25" "  +0.44%  "

# Row 26
Set-TextCell "D26" "0.1230"
Set-TextCell "E26" "  +4.42%  "

# Row 27
Set-TextCell "D27" "15.80"
Set-TextCell "E27" "  +1.02%  "

# Row 28
Set-TextCell "D28" "0.06637"
Set-TextCell "E28" "  +12.97%  "

# Row 29
Set-TextCell "D29" "1.350"
Set-TextCell "E29" "  +6.62%  "

# Row 30
Set-TextCell "D30" "1.329"
Set-TextCell "E30" "  +0.53%  "

# Row 31
Set-TextCell "D31" "3.571"
Set-TextCell "E31" "  +1.63%  "

# Row 32
Set-TextCell "D32" "3.565"
Set-TextCell "E32" "  +0.94%  "

# Row 33
Set-TextCell "D33" "1.657"
Set-TextCell "E33" "  +0.94%  "

# Row 34
Set-TextCell "D34" "1.026"
Set-TextCell "E34" "  +1.29%  "

# Row 35
Set-TextCell "D35" "0.6213"
Set-TextCell "E35" "  +3.50%  "

# Row 36
Set-TextCell "D36" "2.404"
Set-TextCell "E36" "  +2.09%  "

# Row 37
Set-TextCell "D37" "2.697"
Set-TextCell "E37" "  +1.90%  "

# Row 38
Set-TextCell "D38" "6.373"
Set-TextCell "E38" "  +5.94%  "

# Row 39
Set-TextCell "D39" "1.111.17"
Set-TextCell "E39" "  +3.42%  "

# Row 40
Set-TextCell "D40" "0.01616"
Set-TextCell "E40" "  +0.01%  "

# Row 41
Set-TextCell "D41" "0.8806"
Set-TextCell "E41" "  +1.65%  "

# Row 42
Set-TextCell "D42" "1.016"
Set-TextCell "E42" "  +0.69%  "

# Row 43
Set-TextCell "D43" "101.07"
Set-TextCell "E43" "  +1.35%  "

# Row 44
Set-TextCell "D44" "1.833.81"
Set-TextCell "E44" "  +0.88%  "

# Row 45
Set-TextCell "E45" "  +1.08%  "

# Row 46
Set-TextCell "D46" "56.76"
Set-TextCell "E46" "  +1.63%  "

# Row 47
Set-TextCell "D47" "8.166"
Set-TextCell "E47" "  +0.96%  "

# Row 48
Set-TextCell "D48" "1.007"
Set-TextCell "E48" "  +0.34%  "

# Row 49
Set-TextCell "D49" "0.05268"
Set-TextCell "E49" "  +1.69%  "

# Row 51
Set-TextCell "D51" "6.058"
Set-TextCell "E51" "  +3.25%  "
